$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "RN-002: Após o envio do orçamento cliente tem até uma semana para fazer o cancelamento do serviço.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RN-002: Após envio do orçamento ao cliente, ele tem até uma semana para fazer o cancelamento do serviço.",
    2
)

$d.Content.Find.Execute(
    "RN-006: O cliente geralmente aceita na hora ou recusa pela urgência do trabalho, mas tem 48 horas para aceitar.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RN-006: O cliente geralmente aceita de imediato ou recusa pela urgência do trabalho, mas tem até 48 horas para aceitar.",
    2
)

$d.Content.Find.Execute(
    "RN-007: O cliente tem 48 horas para cancelar o serviço após o envio do orçamento.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RN-007: O cliente tem 48 horas para cancelar o serviço após o envio do orçamento pela empresa.",
    2
)
